# Add a new "clientvoucherused" row to the table on the single worksheet,
# and mark a few more tables as having a "view questionnaire" (ok) flag.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 5 ("detailtypes"), shifting
# everything below it down by one. This becomes the new row for
# "clientvoucherused".
$ws.Rows.Item(5).Insert()

# Fill in the new row's table name and "ok" marker.
$ws.Cells.Item(5, 1).Value = "clientvoucherused"
$ws.Cells.Item(5, 2).Value = "ok"

# A few existing rows also gain the "ok" marker in column B.
$ws.Cells.Item(2, 2).Value = "ok"   # clientquestionnaries
$ws.Cells.Item(4, 2).Value = "ok"   # clientvoucher
$ws.Cells.Item(8, 2).Value = "ok"   # logitems

# Update the active selection to match the edited workbook.
$ws.Range("B10").Select()
